$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Force text storage for D so numeric-looking values are not reinterpreted as numbers.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "36.589.00"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -0.41%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.968.21"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  +0.08%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "243.81"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +0.04%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.626"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +1.92%  "

$ws.Range("E7").Value = "  +2.73%  "

$ws.Range("E8").Value = "  +0.00%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.379"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +1.98%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0788"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("E11").Value = "  +0.72%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.23"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +4.43%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.842"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +2.64%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.260.01"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +0.68%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.56"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("E16").Value = "  +0.34%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.970.62"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +0.49%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "36.508.61"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -0.46%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "69.78"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("E20").Value = "  -0.76%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "229.48"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("E23").Value = "  +0.17%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.43"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("E25").Value = "  +1.09%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.145"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +5.65%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.13"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("E28").Value = "  +0.85%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.39"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +0.11%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.34"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +20.85%  "

$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("E33").Value = "  -0.73%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.52"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +6.76%  "

$ws.Range("E35").Value = "  +3.76%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  -2.06%  "

$ws.Range("E38").Value = "  +0.33%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.39"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -13.81%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0964"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -3.59%  "

$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("E43").Value = "  -1.25%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.84"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -1.42%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.367.37"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +1.69%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "88.90"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +1.68%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.02"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -1.36%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.18"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  -0.47%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "45.95"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +5.65%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.153.35"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +0.77%  "
